# Customer_List.xlsx
#
# The "Phone Number" column (D) was re-entered as plain numbers instead of
# the previously hand-formatted text strings like "(562)123-4567" - Excel
# had stored those as shared-string text, but here they become native
# numeric values (e.g. 1621234567). Re-typing the column this way also
# drops the five old formatted-phone shared strings from the string table
# automatically, since nothing references them any more.
#
# Finally, the cursor is left parked on D6 - the last phone-number cell
# that was edited - instead of the old selection down the (now untouched)
# column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new numeric phone value for column D
$phoneNumbers = @{
    2 = 1621234567
    3 = 1235150708
    4 = 1149876543
    5 = 1621281111
    6 = 1261081215
}

foreach ($row in $phoneNumbers.Keys) {
    $ws.Cells.Item($row, 4).Value = $phoneNumbers[$row]
}

$ws.Range("D6").Select()
